# Update the cached "datetimeFigureOut" date placeholder text from
# 4/12/2024 to 4/28/2024 everywhere it appears: the Slide Master and
# every slide layout (CustomLayout) hanging off it.

$p = $ppt.ActivePresentation

$oldDate = "4/12/2024"
$newDate = "4/28/2024"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1. Slide Master
Update-DatePlaceholders $p.SlideMaster.Shapes

# 2. Every slide layout hanging off the Slide Master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholders $layouts.Item($L).Shapes
}
